$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 1")

# Insert a new "Wednesday" row between Tuesday (row 5) and the old Thursday row (row 6),
# pushing Thursday/Friday/Saturday/Sunday (and the trailing note) down by one row.
$ws.Rows.Item(6).Insert()

# Day label for the newly inserted row.
$ws.Range("A6").Value = "Wednesday"

# Re-point the two header cells that used to mirror Totals!D2 / Totals!F2 with
# literal names for the two new team members.
$ws.Range("D2").Value = "Jeff Walters"
$ws.Range("F2").Value = "Gino Betetta"

# Running-total formulas for the new row (mirrors the pattern used by every other row).
$ws.Range("C6").Formula = "=B6+C5"
$ws.Range("E6").Formula = "=D6+E5"
$ws.Range("G6").Formula = "=F6+G5"
$ws.Range("I6").Formula = "=H6+I5"

# Re-anchor the running-total formulas for the rows that shifted down so each one
# references the row immediately above it (the row-insert does not always rewrite
# the shared-formula member cells correctly).
$ws.Range("C7").Formula = "=B7+C6"
$ws.Range("E7").Formula = "=D7+E6"
$ws.Range("G7").Formula = "=F7+G6"
$ws.Range("I7").Formula = "=H7+I6"

$ws.Range("C8").Formula = "=B8+C7"
$ws.Range("E8").Formula = "=D8+E7"
$ws.Range("G8").Formula = "=F8+G7"
$ws.Range("I8").Formula = "=H8+I7"

$ws.Range("C9").Formula = "=B9+C8"
$ws.Range("E9").Formula = "=D9+E8"
$ws.Range("G9").Formula = "=F9+G8"
$ws.Range("I9").Formula = "=H9+I8"

$ws.Range("C10").Formula = "=B10+C9"
$ws.Range("E10").Formula = "=D10+E9"
$ws.Range("G10").Formula = "=F10+G9"
$ws.Range("I10").Formula = "=H10+I9"

# Fill in the updated/new daily hours.
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = 2

$ws.Range("D5").Value = 3
$ws.Range("F5").Value = 3

$ws.Range("B6").Value = 2
$ws.Range("D6").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("H6").Value = 4

$ws.Range("B7").Value = 3
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("H7").Value = 2

$ws.Range("D8").Value = 4
$ws.Range("H8").Value = 2

$ws.Range("H9").Value = 0

$ws.Range("H10").Value = 2

# Match the author's final cursor position on the sheet.
$ws.Range("H11").Select() | Out-Null
